# Repull data, push all data, mean calculation
# Update the dSF column (F) values for rows where newly pulled data differs
# from the previously stored (stale) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    3  = -6
    4  = -2
    9  = 0
    10 = 3
    13 = 4
    14 = -3
    17 = -4
    18 = -9
    20 = -7
    21 = -7
    23 = -6
    29 = -7
    30 = -2
    32 = -6
    36 = 1
    38 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
